# Jobsheet 8 - template_barang.xlsx
# Tugas 1 > import user: replace the sample "barang" rows with the new
# product list (kategori_id now numeric, new barang_kode / barang_nama
# values, updated harga_beli / harga_jual) and tidy up the sheet look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row: drop the old bold/small-font style, back to Normal ----
$ws.Range("A1:E1").Style = "Normal"

# ---- kategori_id (numeric category ids) ----
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3

# ---- barang_kode ----
$ws.Range("B2").Value = "M003"
$ws.Range("B3").Value = "M004"
$ws.Range("B4").Value = "B003"
$ws.Range("B5").Value = "B004"
$ws.Range("B6").Value = "H003"

# ---- barang_nama ----
$ws.Range("C2").Value = "Liptint Velvel"
$ws.Range("C3").Value = "BB Creaam"
$ws.Range("C4").Value = "Body Lotion Bengkoang"
$ws.Range("C5").Value = "Shower Gel Bunga"
$ws.Range("C6").Value = "Shampoo Anti Rontok"

# barang_nama column gets an explicit (black, 11pt) font
$ws.Range("C2:C6").Font.Size = 11
$ws.Range("C2:C6").Font.Color = 0

# ---- harga_beli ----
$ws.Range("D2").Value = 35000
$ws.Range("D3").Value = 65000
$ws.Range("D4").Value = 28000
$ws.Range("D5").Value = 45000
$ws.Range("D6").Value = 37000

# ---- harga_jual ----
$ws.Range("E2").Value = 55000
$ws.Range("E3").Value = 80000
$ws.Range("E4").Value = 40000
$ws.Range("E5").Value = 65000
$ws.Range("E6").Value = 50000

# ---- column widths (A:E tweaked slightly, F:H prepared for more columns) ----
$ws.Columns.Item(1).ColumnWidth = 8.4362
$ws.Columns.Item(2).ColumnWidth = 12.7096
$ws.Columns.Item(3).ColumnWidth = 21.7995
$ws.Columns.Item(4).ColumnWidth = 9.983
$ws.Columns.Item(5).ColumnWidth = 9.6198
$ws.Columns.Item(6).ColumnWidth = 9.7995
$ws.Columns.Item(7).ColumnWidth = 10.2565
$ws.Columns.Item(8).ColumnWidth = 10.4362

# ---- view: zoom to 78% and leave the selection on H6 ----
$excel.ActiveWindow.Zoom = 78
$ws.Range("H6").Select() | Out-Null

# ---- printer settings are no longer relevant, drop them ----
$ws.PageSetup.Orientation = 1 | Out-Null
